$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 38267
$ws.Range("I6").Value = 10001
$ws.Range("J6").Value = 52400
$ws.Range("K6").Value = 30003
$ws.Range("L6").Value = 157200
$ws.Range("M6").Value = -29891
$ws.Range("N6").Value = -157424

$ws.Range("H40").Value = 6870.048
$ws.Range("I40").Value = 9681.615
$ws.Range("J40").Value = 2301.25
$ws.Range("K40").Value = 9681.615
$ws.Range("L40").Value = 2301.25
$ws.Range("M40").Value = -9506.615
$ws.Range("N40").Value = -2651.25

$ws.Range("H64").Value = 3596.075
$ws.Range("I64").Value = 3302.8572
$ws.Range("J64").Value = 3753.9614
$ws.Range("K64").Value = 3302.8572
$ws.Range("L64").Value = 3753.9614
$ws.Range("M64").Value = -3054.8572
$ws.Range("N64").Value = -4249.9614

$ws.Range("H67").Value = 3596.075
$ws.Range("I67").Value = 3302.8572
$ws.Range("J67").Value = 3753.9614
$ws.Range("K67").Value = 3302.8572
$ws.Range("L67").Value = 3753.9614
$ws.Range("M67").Value = -2444.8572
$ws.Range("N67").Value = -5469.9614

$ws.Range("H113").Value = 2076.25
$ws.Range("I113").Value = 1744
$ws.Range("K113").Value = 1744
$ws.Range("M113").Value = 1510

$ws.Range("H129").Value = 906.68805
$ws.Range("J129").Value = 927.35236
$ws.Range("L129").Value = 2782.05708
$ws.Range("N129").Value = -12782.05708

$ws.Range("H137").Value = 1712752.4
$ws.Range("I137").Value = 2849890.8
$ws.Range("K137").Value = 8549672.399999999
$ws.Range("M137").Value = -8547122.399999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 343.5
$ws.Range("I5").Value = 343.5
$ws.Range("K5").Value = 343.5
$ws.Range("M5").Value = -231.5

$ws.Range("H32").Value = 11832.038
$ws.Range("I32").Value = 11002.261
$ws.Range("K32").Value = 11002.261
$ws.Range("M32").Value = -10715.261

$ws.Range("H45").Value = 2420.8235
$ws.Range("I45").Value = 2125
$ws.Range("J45").Value = 3382.25
$ws.Range("K45").Value = 2125
$ws.Range("L45").Value = 3382.25
$ws.Range("M45").Value = -1748
$ws.Range("N45").Value = -4136.25

$ws.Range("H88").Value = 17871058
$ws.Range("I88").Value = 28575086
$ws.Range("K88").Value = 28575086
$ws.Range("M88").Value = -28574680

$ws.Range("H91").Value = 17871058
$ws.Range("I91").Value = 28575086
$ws.Range("K91").Value = 28575086
$ws.Range("M91").Value = -28573682

$ws.Range("H110").Value = 1345.9062
$ws.Range("I110").Value = 1289
$ws.Range("K110").Value = 1289
$ws.Range("M110").Value = 756

$ws.Range("H122").Value = 2046.6666
$ws.Range("I122").Value = 2060
$ws.Range("J122").Value = 1980
$ws.Range("K122").Value = 6180
$ws.Range("L122").Value = 5940
$ws.Range("M122").Value = -3730
$ws.Range("N122").Value = -10840

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 343.5
$ws.Range("I4").Value = 343.5
$ws.Range("K4").Value = 343.5
$ws.Range("M4").Value = -228.5

$ws.Range("H86").Value = 1943.8889
$ws.Range("I86").Value = 2066.6667
$ws.Range("J86").Value = 1882.5
$ws.Range("K86").Value = 2066.6667
$ws.Range("L86").Value = 1882.5
$ws.Range("M86").Value = -943.6667000000002
$ws.Range("N86").Value = -4128.5

$ws.Range("H89").Value = 1943.8889
$ws.Range("I89").Value = 2066.6667
$ws.Range("J89").Value = 1882.5
$ws.Range("K89").Value = 10333.3335
$ws.Range("L89").Value = 9412.5
$ws.Range("M89").Value = -4717.333500000001
$ws.Range("N89").Value = -20644.5

$ws.Range("H107").Value = 2046.5454
$ws.Range("I107").Value = 1616.5
$ws.Range("J107").Value = 2562.6
$ws.Range("K107").Value = 1616.5
$ws.Range("L107").Value = 2562.6
$ws.Range("M107").Value = 303.5
$ws.Range("N107").Value = -6402.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1985.32
$ws.Range("I31").Value = 792.8261
$ws.Range("J31").Value = 3001.1482
$ws.Range("K31").Value = 792.8261
$ws.Range("L31").Value = 3001.1482
$ws.Range("M31").Value = -497.8261
$ws.Range("N31").Value = -3591.1482

$ws.Range("H34").Value = 1985.32
$ws.Range("I34").Value = 792.8261
$ws.Range("J34").Value = 3001.1482
$ws.Range("K34").Value = 792.8261
$ws.Range("L34").Value = 3001.1482
$ws.Range("M34").Value = -590.8261
$ws.Range("N34").Value = -3405.1482

$ws.Range("H86").Value = 3971.5454
$ws.Range("I86").Value = 3971.5454
$ws.Range("K86").Value = 3971.5454
$ws.Range("M86").Value = -2848.5454

$ws.Range("H87").Value = 9000
$ws.Range("J87").Value = 9000
$ws.Range("L87").Value = 9000
$ws.Range("N87").Value = -11372

$ws.Range("H89").Value = 3971.5454
$ws.Range("I89").Value = 3971.5454
$ws.Range("K89").Value = 19857.727
$ws.Range("M89").Value = -14241.727

$ws.Range("H90").Value = 9000
$ws.Range("J90").Value = 9000
$ws.Range("L90").Value = 27000
$ws.Range("N90").Value = -38856

$ws.Range("H141").Value = 7948.5
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 7948.5
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 7948.5
$ws.Range("N141").Value = -18308.5
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3581.5144
$ws.Range("I5").Value = 4158.074
$ws.Range("J5").Value = 1635.625
$ws.Range("K5").Value = 12474.222
$ws.Range("L5").Value = 4906.875
$ws.Range("M5").Value = -12362.222
$ws.Range("N5").Value = -5130.875

$ws.Range("H105").Value = 170416.33
$ws.Range("J105").Value = 203999.6
$ws.Range("L105").Value = 611998.8
$ws.Range("N105").Value = -617240.8

$ws.Range("H135").Value = 3581.5144
$ws.Range("I135").Value = 4158.074
$ws.Range("J135").Value = 1635.625
$ws.Range("K135").Value = 37422.666
$ws.Range("L135").Value = 14720.625
$ws.Range("M135").Value = -34887.666
$ws.Range("N135").Value = -19790.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1423.1428
$ws.Range("I122").Value = 1103.4286
$ws.Range("K122").Value = 3310.2858
$ws.Range("M122").Value = -860.2857999999997

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 113334.336
$ws.Range("I122").Value = 113334.336
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 340003.008
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -337553.008
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5650
$ws.Range("I81").Value = 1122.2222
$ws.Range("J81").Value = 13800
$ws.Range("K81").Value = 2244.4444
$ws.Range("L81").Value = 27600
$ws.Range("M81").Value = -1183.4444
$ws.Range("N81").Value = -29722

$ws.Range("H84").Value = 5650
$ws.Range("I84").Value = 1122.2222
$ws.Range("J84").Value = 13800
$ws.Range("K84").Value = 11222.222
$ws.Range("L84").Value = 138000
$ws.Range("M84").Value = -5918.222
$ws.Range("N84").Value = -148608
